{"js": "const body = context.document.body;\n\n// 1) Title: \"Evidence for recovery of an ant-dispersed myrmecochorous plant community\n//    following a small-scale disturbance\"\n//    -> \"An ant-dispersed plant community recovers following a small-scale disturbance\n//        in a Connecticut forest\"\nconst titleResults = body.search(\n  \"Evidence for recovery of an ant-dispersed myrmecochorous plant community following a small-scale disturbance\",\n  { matchCase: true }\n);\ntitleResults.load(\"text\");\n\n// 2) Abstract sentence: \"nutrient rich appendages are used to recruit ants...\"\n//    -> \"nutrient-rich seed appendages recruit ants...\"\nconst abstractResults = body.search(\n  \"Many species of plants exhibit a dispersal syndrome called myrmecochory, in which nutrient rich appendages are used to recruit ants which in turn transport seeds.\",\n  { matchCase: true }\n);\nabstractResults.load(\"text\");\n\n// 3) \"6 year duration\" -> \"6-year duration\"\nconst durationResults = body.search(\n  \"had made a recovery within the 6 year duration of the experiment.\",\n  { matchCase: true }\n);\ndurationResults.load(\"text\");\n\n// 4) \"surveys were to provide\" -> \"surveys was to provide\"\nconst surveysResults = body.search(\n  \"The goal of these surveys were to provide preliminary data\",\n  { matchCase: true }\n);\nsurveysResults.load(\"text\");\n\nawait context.sync();\n\ntitleResults.items[0].insertText(\n  \"An ant-dispersed plant community recovers following a small-scale disturbance in a Connecticut forest\",\n  \"Replace\"\n);\n\nabstractResults.items[0].insertText(\n  \"Many species of plants exhibit a dispersal syndrome called myrmecochory, in which nutrient-rich seed appendages recruit ants which in turn transport seeds.\",\n  \"Replace\"\n);\n\ndurationResults.items[0].insertText(\n  \"had made a recovery within the 6-year duration of the experiment.\",\n  \"Replace\"\n);\n\nsurveysResults.items[0].insertText(\n  \"The goal of these surveys was to provide preliminary data\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# 1) Title: \"Evidence for recovery of an ant-dispersed myrmecochorous plant community\n#    following a small-scale disturbance\"\n#    -> \"An ant-dispersed plant community recovers following a small-scale disturbance\n#        in a Connecticut forest\"\nReplace-Text \"Evidence for recovery of an ant-dispersed myrmecochorous plant community following a small-scale disturbance\" \"An ant-dispersed plant community recovers following a small-scale disturbance in a Connecticut forest\"\n\n# 2) Abstract sentence: \"nutrient rich appendages are used to recruit ants...\"\n#    -> \"nutrient-rich seed appendages recruit ants...\"\nReplace-Text \"Many species of plants exhibit a dispersal syndrome called myrmecochory, in which nutrient rich appendages are used to recruit ants which in turn transport seeds.\" \"Many species of plants exhibit a dispersal syndrome called myrmecochory, in which nutrient-rich seed appendages recruit ants which in turn transport seeds.\"\n\n# 3) \"6 year duration\" -> \"6-year duration\"\nReplace-Text \"had made a recovery within the 6 year duration of the experiment.\" \"had made a recovery within the 6-year duration of the experiment.\"\n\n# 4) \"surveys were to provide\" -> \"surveys was to provide\"\nReplace-Text \"The goal of these surveys were to provide preliminary data\" \"The goal of these surveys was to provide preliminary data\"\n"}
